# Reorder the six scene-category columns (A:F) on the active sheet.
#
# The sheet holds a 1-hot "block order" table: row 1 is a header naming
# each scene-category slot (e.g. "kitchens_1"), rows 2-7 hold the 0/1
# assignment for each subject/row. The edit re-orders the *columns*
# (header + all data underneath it move together) into this new order:
#   kitchens_1, kitchens_2, living_rooms_1, bedrooms_1, bedrooms_2, living_rooms_2
#
# Old column -> New column mapping (derived from the column contents):
#   old A (living_rooms_1) -> new C
#   old B (kitchens_1)     -> new A
#   old C (bedrooms_1)     -> new D
#   old D (kitchens_2)     -> new B
#   old E (living_rooms_2) -> new F
#   old F (bedrooms_2)     -> new E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldOrder = @("A", "B", "C", "D", "E", "F")
$newOrder = @("B", "D", "A", "C", "F", "E")  # new col i comes from old column $newOrder[i]

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 7) { $lastRow = 7 }

# Snapshot the full used range (header + data) before writing anything,
# since source and destination columns overlap. Use .Value2 for reads
# (.Value's getter surfaces a property descriptor, not the cell content,
# in this host) and .Value for writes.
$snapshot = @{}
foreach ($col in $oldOrder) {
    $colValues = @()
    for ($row = 1; $row -le $lastRow; $row++) {
        $colValues += , ($ws.Range("$col$row").Value2)
    }
    $snapshot[$col] = $colValues
}

# Write the snapshot back out in the new column order.
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $destCol = $oldOrder[$i]
    $srcCol = $newOrder[$i]
    $values = $snapshot[$srcCol]
    for ($row = 1; $row -le $lastRow; $row++) {
        $ws.Range("$destCol$row").Value = $values[$row - 1]
    }
}
